$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EEU data")

# Update car uncomp elasticity (column K, "e_qs_ps_UC") from -0.1 to -0.2
# for the two "Car" rows (rows 2 and 3), per PB change.
$ws.Range("K2").Value = -0.2
$ws.Range("K3").Value = -0.2

# Reflect the resulting active cell/selection on the sheet
[void]$ws.Range("K4").Select()
